$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 new rows before row 2 for the new grant entry (2023 grant,
#    "Efecto del control de los recursos real y simulado..."), pushing all
#    existing grant rows down by 3.
# ---------------------------------------------------------------------------
$ws.Rows("2:4").Insert()

# New row 2: the grant header row (what / when / with / where / why)
$ws.Range("A2").Value = 'XI \href{https://www.unbosque.edu.co/investigaciones/convocatorias-investigacion}{Convocatoria Interna para la Financiación de Proyectos de Investigación}, 2023'
$ws.Range("B2").Value = 'Feb. 2024 - Actualmente'
$ws.Range("C2").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}'
$ws.Range("D2").Value = 'Bogota, Colombia'
$ws.Range("E2").Value = 'Proyecto: \textit{Efecto del control de los recursos real y simulado sobre las preferencias de mujeres andrófilas por la masculinidad en rostros de hombres: un estudio experimental usando rastreo ocular}'

# New row 3: principal investigator note (plain wrap style, same as row 2)
$ws.Range("E3").Value = 'Investigadora principal: \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}'

# New row 4: grant amount (currency style)
$ws.Range("E4").Value = 'COP\$89.979.750'
$ws.Range("E4").NumberFormat = "_-[`$`$-240A]\ * #,##0.00_-;\-[`$`$-240A]\ * #,##0.00_-;_-[`$`$-240A]\ * ""-""??_-;_-@_-"
$ws.Range("E4").HorizontalAlignment = -4131
$ws.Range("E4").VerticalAlignment = -4160
$ws.Range("E4").WrapText = $true

# Row height for the new grant header row (matches the other 60pt grant rows)
$ws.Rows("2").RowHeight = 60

# ---------------------------------------------------------------------------
# 2. Wrap the existing "Proyecto: ..." descriptions with \textit{...} for the
#    three older grants (each description appears twice in the sheet).
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = 'Proyecto: \textit{Señales perceptibles de salud física y mental en rostros, voces y olores corporales, y su relación con niveles hormonales}'
$ws.Range("E7").Value = 'Proyecto: \textit{Señales perceptibles de salud física y mental en rostros, voces y olores corporales, y su relación con niveles hormonales}'
$ws.Range("E9").Value = 'Proyecto: \textit{Efecto de señales estáticas evolutivamente relevantes (sexo, dominancia y atractivo) en el procesamiento cortical de rostros humanos}'
$ws.Range("E11").Value = 'Proyecto: \textit{Efectos de los niveles hormonales, masculinidad y feminidad, en la discriminación tonal en hombres y mujeres}'

# ---------------------------------------------------------------------------
# 3. Widen column B to fit the new, longer "when" values.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 35.15

# ---------------------------------------------------------------------------
# 4. Update selection to match the saved view state.
# ---------------------------------------------------------------------------
[void]$ws.Range("B7").Select()
